$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# POL-3138: Stricter Hub finding.
# Origin/Destination country + city + locode are now always populated together
# (previously some rows only had the locode, others only had the country/city).

$rows = @(2, 3, 4)

foreach ($r in $rows) {
    $ws.Range("E$r").Value2 = "Sweden"
    $ws.Range("F$r").Value2 = "Gothenburg"
    $ws.Range("G$r").Value2 = "SEGOT"
    $ws.Range("H$r").Value2 = "China"
    $ws.Range("I$r").Value2 = "Shanghai"
    $ws.Range("J$r").Value2 = "CNSHA"
    $ws.Range("J$r").Font.Name = "Arial"
}

# Transshipment locode is only present on row 2
$ws.Range("K2").Value2 = "ZACPT"

# Column widths shrink slightly after the data change
$ws.Range("E2:G4").ColumnWidth = 15.55
$ws.Range("H2:J4").ColumnWidth = 21.3
